# Update Name of Algo
# Apply targeted numeric corrections to the result_data_RandomForest sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value  = 12.46509999999999
$ws.Range("E7").Value  = 12.1094
$ws.Range("C8").Value  = -11.10889999999999
$ws.Range("B12").Value = 5.687699999999998
$ws.Range("C12").Value = -14.68270000000002
$ws.Range("C14").Value = -11.94499999999999
$ws.Range("E19").Value = 13.1435
$ws.Range("E21").Value = 12.77329999999999
$ws.Range("C22").Value = -10.66739999999998
$ws.Range("E24").Value = 12.84859999999999
